# Generate Report for Handoff
# Update status text and timestamps, and shrink status-datetime column widths.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# Status text (shared with zh-cn/de-de "Status" column): "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-08-22 01:09:14"

# Shrink the zh-cn / de-de status columns on Overview
$wsOverview.Range("E1").ColumnWidth = 17.2159881591797
$wsOverview.Range("F1").ColumnWidth = 17.2159881591797

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-22 01:09:10"
$wsZhCn.Range("C1").ColumnWidth = 17.2159881591797

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-22 01:09:14"
$wsDeDe.Range("C1").ColumnWidth = 17.2159881591797
